$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date format (style) from A52 down to the newly added A53 row before filling values
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 3.14593994906931
$ws.Range("D2").Value = 2008

$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("D3").Value = 2009

$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 1.769627576887389
$ws.Range("D4").Value = 2009

$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("D5").Value = 2010

$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = -4.774178217057756
$ws.Range("D6").Value = 2010

$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 0.26989196962941
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 0.5668332592311476

$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = 1.97975191822708
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 2.743551941645217

$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 3.872616460773104
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 2.302293339151018

$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 3.452886745653183
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 1.794132456841213

$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 1.385102141502959
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 1.418650881502459

$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 1.239479831392853
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 1.158731032337301

$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = -0.09990319152052285
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 0.9503446311848185

$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = 0.2379616621361214
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 1.348985046565354

$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 1.985365612881851
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 1.683270396159919

$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 1.51977456621637
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 0.7739869831243862

$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 1.326941625882871
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 1.183026708734536

$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 1.470039379455756
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 1.319057785023592

$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = 1.606236217798274
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 1.5164178751083

$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = 1.638797242243251
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 1.236938064849924

$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 1.730343268967593
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 1.500360889159746

$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 2.161565493242668
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 2.320541194291881

$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 2.353897533252747
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 1.486421985254816

$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 2.175463816693268
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 1.290682042302871

$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 2.306979482028937
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 1.375814392482444

$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 2.214251681313772
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 0.6374750548026054

$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 0.4816924525496313
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 0.9776390272966617

$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 0.7385331577992593
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 1.193318741914795

$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = 0.6388168203515399
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = 0.7331704885201074

$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = 0.6066442151010376
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 0.57214245765278

$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 0.4580828524435532
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 0.9346503642840398

$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -1.986210268830169
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = -1.030455917249229

$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = -4.65090747647452
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = -2.124426335989094

$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = -4.207901339433196
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = -0.2586890779524231

$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = 0.590529387845784
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = 0.9817808318997479

$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = 0.4846423081591666
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = 1.007198277338284

$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = 0.7583924418458787
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 1.422024064575078

$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 1.099928004397532
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 1.794400784768979

$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 2.197506404324789
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = 1.345036064408078

$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 1.990013243928312
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 1.178997008351645

$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = 2.236860175919531
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 1.139002134681211

$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = 2.310042359896225
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = 0.8232644777432796

$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = 0.3858398712458078
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 0.9487085336255197

$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = -0.07210020592836042
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 0.7849798646630823

$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = -0.09588622947416248
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 0.5889550719078596

$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = 0.0464415346324687
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = 0.4457784880425031

$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = -0.03414634355979329
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = 0.8055775472706417

$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = -0.03189435474734159
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 0.8280060478212947

$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = -0.2385784141923808
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = 0.6469810241417351

$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = -0.3101476031197148
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = 0.5215192790195111

$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = 0.16526024324508
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = 0.8471324137523606

$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 0.1720146172997206
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = 0.822899916144304

$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = -0.0960403240804597
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 0.6335680123873866

# C5 no longer has a value in the updated dataset
$ws.Range("C5").Value = ""
